$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("H2").Value = 3.9
$ws.Range("I2").Value = 3.3
$ws.Range("N2").Value = 17
$ws.Range("AB2").Value = 19
$ws.Range("AI2").Value = 21
